$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column G to fit the new comment text
$ws.Columns.Item(7).ColumnWidth = 30.3

# Row 12 - was blank, now filled in with "Shortest Common Supersequence"
$ws.Range("A12").Value = 11
$ws.Range("B12").Value = "Shortest Common Supersequence"
$ws.Range("C12").Value = "DP"
$ws.Range("D12").Value = "String"
$ws.Range("E12").Value = "medium"
$ws.Range("F12").Value = "GeeksForGeeks"

# Row 13 - new row, with highlighted comment in G13
$ws.Range("A13").Value = 12
$ws.Range("B13").Value = "Maximum number of A’s using given four keys"
$ws.Range("C13").Value = "DP"
$ws.Range("D13").Value = "Math"
$ws.Range("E13").Value = "medium"
$ws.Range("F13").Value = "GeeksForGeeks"
$ws.Range("G13").Value = "think about the transform func"

# New row's cells inherit style (centered) like rest of table, plus G13 highlighted yellow
$ws.Range("A12:F13").HorizontalAlignment = -4108
$ws.Range("G13").HorizontalAlignment = -4108
$ws.Range("G13").Interior.Color = 65535

# Update selection to match the new active cell noted in the diff
$ws.Range("E16").Select()
